$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the 1-based index of the paragraph whose text contains $Needle.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex {
    param([string]$Needle)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.Contains($Needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Helper: replace the text between character offsets [$LocalStart,$LocalEnd)
# of the paragraph at $ParaIndex (offsets local to that paragraph's own
# Range.Text) with $NewSub, leaving the untouched text before/after it as
# separate runs (matching how Word itself splits a run when you edit in the
# middle of it) rather than collapsing the whole paragraph back into a
# single run.
#
# We build the replacement text in a throw-away paragraph inserted right
# after the target paragraph (so it inherits the *target's* clean run
# formatting, not whatever trailing formatting happens to sit at the very
# end of the document), copy that text's FormattedText into the exact
# sub-range being replaced (which keeps run boundaries at the splice
# points), then delete the scratch paragraph again.
# ---------------------------------------------------------------------------
function Split-ReplaceInParagraph {
    param([int]$ParaIndex, [int]$LocalStart, [int]$LocalEnd, [string]$NewSub)

    $target = $d.Paragraphs($ParaIndex).Range
    $paraStart = $target.Start
    $replaceRange = $d.Range($paraStart + $LocalStart, $paraStart + $LocalEnd)

    # Scratch paragraph right after the target (clean, matching formatting).
    $target.InsertParagraphAfter()
    $scratch = $d.Paragraphs($ParaIndex + 1).Range
    $scratch.InsertAfter($NewSub)
    # A paragraph's Range (and a range grown via InsertAfter from an empty
    # paragraph) includes the trailing paragraph-mark character; trim it off
    # so we don't splice a hidden pilcrow into the target paragraph.
    $scratch.MoveEnd(1, -1) | Out-Null

    # Splice the scratch run(s) into the target sub-range; this keeps the
    # surrounding text of $target as separate, untouched runs.
    $replaceRange.FormattedText = $scratch.FormattedText

    # Clean up the scratch paragraph (still immediately follows the edited
    # target paragraph).
    $d.Paragraphs($ParaIndex + 1).Range.Delete()
}

# ---------------------------------------------------------------------------
# "Location: Geographic location, ..." -> "Location: Geographical location, ..."
# i.e. insert "al" right after "...Geographic" (before " location, ...").
# ---------------------------------------------------------------------------
$p1 = Find-ParagraphIndex "Location: Geographic location, typically the country name followed by a name of the region"
$p1Text = $d.Paragraphs($p1).Range.Text
$splitAt = $p1Text.IndexOf("Location: Geographic") + "Location: Geographic".Length
Split-ReplaceInParagraph $p1 $splitAt $splitAt "al"

# ---------------------------------------------------------------------------
# "in which the population was collected" -> "in which the insect was collected"
# i.e. replace "population" with "insect".
# ---------------------------------------------------------------------------
$p2 = Find-ParagraphIndex "in which the population was collected"
$p2Text = $d.Paragraphs($p2).Range.Text
$oldSub = "population"
$start2 = $p2Text.IndexOf($oldSub)
$end2 = $start2 + $oldSub.Length
Split-ReplaceInParagraph $p2 $start2 $end2 "insect"
